$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 461.25
$ws.Range("I12").Value = 392
$ws.Range("J12").Value = 530.5
$ws.Range("K12").Value = 392
$ws.Range("L12").Value = 530.5
$ws.Range("M12").Value = -222
$ws.Range("N12").Value = -870.5
$ws.Range("H15").Value = 1623.7317
$ws.Range("I15").Value = 1623.7317
$ws.Range("K15").Value = 4871.1951
$ws.Range("M15").Value = -4702.1951
$ws.Range("H19").Value = 524.72
$ws.Range("I19").Value = 374
$ws.Range("J19").Value = 663.8461
$ws.Range("K19").Value = 374
$ws.Range("L19").Value = 663.8461
$ws.Range("M19").Value = -199
$ws.Range("N19").Value = -1013.8461
$ws.Range("H29").Value = 786
$ws.Range("I29").Value = 176.18182
$ws.Range("J29").Value = 4140
$ws.Range("K29").Value = 528.5454599999999
$ws.Range("L29").Value = 12420
$ws.Range("M29").Value = -247.5454599999999
$ws.Range("N29").Value = -12982
$ws.Range("H116").Value = 3558.2856
$ws.Range("I116").Value = 3609.1667
$ws.Range("J116").Value = 3253
$ws.Range("K116").Value = 3609.1667
$ws.Range("L116").Value = 3253
$ws.Range("M116").Value = -167.1667000000002
$ws.Range("N116").Value = -10137
$ws.Range("H125").Value = 2236.4443
$ws.Range("I125").Value = 841.6667
$ws.Range("J125").Value = 5026
$ws.Range("K125").Value = 7575.0003
$ws.Range("L125").Value = 45234
$ws.Range("M125").Value = -5115.0003
$ws.Range("N125").Value = -50154
$ws.Range("H137").Value = 9617103
$ws.Range("I137").Value = 1449.7028
$ws.Range("J137").Value = 33335714
$ws.Range("K137").Value = 4349.1084
$ws.Range("L137").Value = 100007142
$ws.Range("M137").Value = -1799.1084
$ws.Range("N137").Value = -100012242
$ws.Range("H138").Value = 2737.47
$ws.Range("I138").Value = 1135.8485
$ws.Range("J138").Value = 3526.3284
$ws.Range("K138").Value = 3407.5455
$ws.Range("L138").Value = 10578.9852
$ws.Range("M138").Value = 1732.4545
$ws.Range("N138").Value = -20858.9852

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H134").Value = 59800
$ws.Range("J134").Value = 59800
$ws.Range("L134").Value = 59800
$ws.Range("N134").Value = -69940

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 7952.75
$ws.Range("J26").Value = 5000
$ws.Range("L26").Value = 5000
$ws.Range("N26").Value = -5584
$ws.Range("H86").Value = 1996.871
$ws.Range("I86").Value = 1835.0526
$ws.Range("J86").Value = 2253.0833
$ws.Range("K86").Value = 1835.0526
$ws.Range("L86").Value = 2253.0833
$ws.Range("M86").Value = -712.0526
$ws.Range("N86").Value = -4499.0833
$ws.Range("H89").Value = 1996.871
$ws.Range("I89").Value = 1835.0526
$ws.Range("J89").Value = 2253.0833
$ws.Range("K89").Value = 9175.262999999999
$ws.Range("L89").Value = 11265.4165
$ws.Range("M89").Value = -3559.262999999999
$ws.Range("N89").Value = -22497.4165
$ws.Range("H94").Value = 1727.2632
$ws.Range("I94").Value = 1741.3334
$ws.Range("J94").Value = 1674.5
$ws.Range("K94").Value = 1741.3334
$ws.Range("L94").Value = 1674.5
$ws.Range("M94").Value = -1290.3334
$ws.Range("N94").Value = -2576.5
$ws.Range("H96").Value = 8037.6
$ws.Range("I96").Value = 2679.25
$ws.Range("K96").Value = 2679.25
$ws.Range("M96").Value = 66.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2166.25
$ws.Range("I99").Value = 1904.5
$ws.Range("J99").Value = 3475
$ws.Range("K99").Value = 1904.5
$ws.Range("L99").Value = 3475
$ws.Range("M99").Value = -406.5
$ws.Range("N99").Value = -6471
$ws.Range("H126").Value = 2166.25
$ws.Range("I126").Value = 1904.5
$ws.Range("J126").Value = 3475
$ws.Range("K126").Value = 5713.5
$ws.Range("L126").Value = 10425
$ws.Range("M126").Value = -3243.5
$ws.Range("N126").Value = -15365
$ws.Range("H140").Value = 52272.5
$ws.Range("J140").Value = 52272.5
$ws.Range("L140").Value = 52272.5
$ws.Range("N140").Value = -62632.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 66724.336
$ws.Range("I11").Value = 80049.2
$ws.Range("K11").Value = 240147.6
$ws.Range("M11").Value = -240007.6
$ws.Range("H22").Value = 33333332
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 33333332
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H68").Value = 1695.6562
$ws.Range("I68").Value = 1636.4419
$ws.Range("J68").Value = 1743.6981
$ws.Range("K68").Value = 4909.3257
$ws.Range("L68").Value = 5231.094300000001
$ws.Range("M68").Value = -4098.3257
$ws.Range("N68").Value = -6853.094300000001
$ws.Range("H71").Value = 1695.6562
$ws.Range("I71").Value = 1636.4419
$ws.Range("J71").Value = 1743.6981
$ws.Range("K71").Value = 14727.9771
$ws.Range("L71").Value = 15693.2829
$ws.Range("M71").Value = -10671.9771
$ws.Range("N71").Value = -23805.2829
$ws.Range("H131").Value = 9636573
$ws.Range("I131").Value = 29471554
$ws.Range("J131").Value = 2438.8286
$ws.Range("K131").Value = 88414662
$ws.Range("L131").Value = 7316.485799999999
$ws.Range("M131").Value = -88409622
$ws.Range("N131").Value = -17396.4858
$ws.Range("H132").Value = 792.8570999999999
$ws.Range("I132").Value = 650
$ws.Range("J132").Value = 900
$ws.Range("K132").Value = 5850
$ws.Range("L132").Value = 8100
$ws.Range("M132").Value = -3320
$ws.Range("N132").Value = -13160

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").ClearContents()
$ws.Range("N75").Value = 0
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").ClearContents()
$ws.Range("N78").Value = 0
$ws.Range("H122").Value = 1182845.6
$ws.Range("J122").Value = 2250
$ws.Range("L122").Value = 6750
$ws.Range("N122").Value = -11650
$ws.Range("H138").Value = 27700.75
$ws.Range("J138").Value = 27700.75
$ws.Range("L138").Value = 27700.75
$ws.Range("N138").Value = -37980.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 274.83334
$ws.Range("I22").Value = 249.75
$ws.Range("J22").Value = 325
$ws.Range("K22").Value = 249.75
$ws.Range("L22").Value = 325
$ws.Range("M22").Value = 45.25
$ws.Range("N22").Value = -915
$ws.Range("H27").Value = 274.83334
$ws.Range("I27").Value = 249.75
$ws.Range("J27").Value = 325
$ws.Range("K27").Value = 249.75
$ws.Range("L27").Value = 325
$ws.Range("M27").Value = -142.75
$ws.Range("N27").Value = -539
$ws.Range("H29").Value = 10488.333
$ws.Range("I29").Value = 8016
$ws.Range("J29").Value = 11724.5
$ws.Range("K29").Value = 8016
$ws.Range("L29").Value = 11724.5
$ws.Range("M29").Value = -7721
$ws.Range("N29").Value = -12314.5
$ws.Range("H46").Value = 379.44446
$ws.Range("I46").Value = 323.33334
$ws.Range("J46").Value = 491.66666
$ws.Range("K46").Value = 323.33334
$ws.Range("L46").Value = 491.66666
$ws.Range("M46").Value = -135.33334
$ws.Range("N46").Value = -867.66666
$ws.Range("H93").Value = 1441.9412
$ws.Range("I93").Value = 896.2857
$ws.Range("J93").Value = 2323.3845
$ws.Range("K93").Value = 896.2857
$ws.Range("L93").Value = 2323.3845
$ws.Range("M93").Value = 351.7143
$ws.Range("N93").Value = -4819.3845
$ws.Range("H132").Value = 1940616.1
$ws.Range("I132").Value = 2606643.8
$ws.Range("J132").Value = 3081.6365
$ws.Range("K132").Value = 7819931.399999999
$ws.Range("L132").Value = 9244.9095
$ws.Range("M132").Value = -7817401.399999999
$ws.Range("N132").Value = -14304.9095

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1286.875
$ws.Range("I81").Value = 1326.3636
$ws.Range("J81").Value = 1200
$ws.Range("K81").Value = 2652.7272
$ws.Range("L81").Value = 2400
$ws.Range("M81").Value = -1591.7272
$ws.Range("N81").Value = -4522
$ws.Range("H84").Value = 1286.875
$ws.Range("I84").Value = 1326.3636
$ws.Range("J84").Value = 1200
$ws.Range("K84").Value = 13263.636
$ws.Range("L84").Value = 12000
$ws.Range("M84").Value = -7959.635999999999
$ws.Range("N84").Value = -22608
$ws.Range("H136").Value = 2629.8958
$ws.Range("I136").Value = 2974.6428
$ws.Range("J136").Value = 2147.25
$ws.Range("K136").Value = 8923.928400000001
$ws.Range("L136").Value = 6441.75
$ws.Range("M136").Value = -6373.928400000001
$ws.Range("N136").Value = -11541.75

Write-Output "Applied all cell updates"
